# Insert a new weekly price record for "Cebollín" (Femacal de La Calera) at
# row 500, pushing the existing rows 500:622 down to 501:623 (dimension
# grows from A1:R622 to A1:R623).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 500..622 down by one, creating a blank row 500.
$ws.Rows.Item(500).Insert()

# Populate the newly inserted row 500 with the new record. Columns that are
# constant for every row of this sheet/subset (market, region, product,
# unit, origin, package size, classification) are repeated as-is; the
# record-specific columns (date, quality, volume, min/max/avg price, $/Kg)
# hold the new observation's values.
$ws.Range("A500").Value = 3
$ws.Range("B500").Value = "Femacal de La Calera"
$ws.Range("C500").Value = "Coquimbo"
$ws.Range("D500").Value = 44943
$ws.Range("E500").Value = 5
$ws.Range("F500").Value = 100112037
$ws.Range("G500").Value = "Cebollín"
$ws.Range("H500").Value = "Sin especificar"
$ws.Range("I500").Value = "Primera"
$ws.Range("J500").Value = 280
$ws.Range("K500").Value = 3300
$ws.Range("L500").Value = 3500
$ws.Range("M500").Value = 3393
$ws.Range("N500").Value = "`$/paquete 36 unidades"
$ws.Range("O500").Value = "Provincia de Quillota"
$ws.Range("P500").Value = 94
$ws.Range("Q500").Value = 36
$ws.Range("R500").Value = "Hortaliza"
